# Add hyperlinks to the three "https://bit.ly/..." sponsor links on the
# "From our sponsor" slide (slide 3). Each of those text runs is the last
# paragraph in its rectangle shape's text frame; PowerPoint represents a
# hyperlink on a run via an ActionSettings/Hyperlink object whose Address
# equals the link target - this also emits the trailing <a:endParaRPr>
# that PowerPoint adds to the paragraph once it has run-level formatting.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$links = @(
    @{ Shape = 2; Paragraph = 3; Url = "https://bit.ly/2IPOjr8" },
    @{ Shape = 4; Paragraph = 3; Url = "https://bit.ly/2FZSKSB" },
    @{ Shape = 6; Paragraph = 3; Url = "https://bit.ly/2GiRGJ4" }
)

foreach ($link in $links) {
    $shape = $s.Shapes.Item($link.Shape)
    $textRange = $shape.TextFrame.TextRange
    $paragraph = $textRange.Paragraphs($link.Paragraph)
    $actionSetting = $paragraph.ActionSettings.Item(1)
    $actionSetting.Hyperlink.Address = $link.Url
}
